$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Price/Volume columns to text format so numeric-looking strings
# (e.g. "4.00", "0.170") keep their exact textual representation.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '68.872.60'  # Price
$ws.Range("E2").Value = '  +1.73%  '  # Volume(1h)
$ws.Range("D3").Value = '3.868.63'  # Price
$ws.Range("E3").Value = '  +0.92%  '  # Volume(1h)
$ws.Range("E4").Value = '  -0.10%  '  # Volume(1h)
$ws.Range("D5").Value = '602.68'  # Price
$ws.Range("E5").Value = '  +1.22%  '  # Volume(1h)
$ws.Range("D6").Value = '171.97'  # Price
$ws.Range("E6").Value = '  +3.88%  '  # Volume(1h)
$ws.Range("D7").Value = '3.866.70'  # Price
$ws.Range("E7").Value = '  +0.90%  '  # Volume(1h)
$ws.Range("E8").Value = '  +0.00%  '  # Volume(1h)
$ws.Range("E9").Value = '  +1.30%  '  # Volume(1h)
$ws.Range("D10").Value = '0.170'  # Price
$ws.Range("E10").Value = '  +4.02%  '  # Volume(1h)
$ws.Range("D11").Value = '6.51'  # Price
$ws.Range("E11").Value = '  +4.12%  '  # Volume(1h)
$ws.Range("E12").Value = '  +1.96%  '  # Volume(1h)
$ws.Range("D13").Value = '0.0000289'  # Price
$ws.Range("E13").Value = '  +17.49%  '  # Volume(1h)
$ws.Range("D14").Value = '37.32'  # Price
$ws.Range("E14").Value = '  +1.88%  '  # Volume(1h)
$ws.Range("D15").Value = '4.522.52'  # Price
$ws.Range("E15").Value = '  +0.99%  '  # Volume(1h)
$ws.Range("D16").Value = '3.850.41'  # Price
$ws.Range("E16").Value = '  +0.27%  '  # Volume(1h)
$ws.Range("D17").Value = '68.837.78'  # Price
$ws.Range("E17").Value = '  +1.64%  '  # Volume(1h)
$ws.Range("D18").Value = '18.51'  # Price
$ws.Range("E18").Value = '  +1.77%  '  # Volume(1h)
$ws.Range("D19").Value = '7.43'  # Price
$ws.Range("E19").Value = '  +0.39%  '  # Volume(1h)
$ws.Range("E20").Value = '  +0.62%  '  # Volume(1h)
$ws.Range("D21").Value = '11.17'  # Price
$ws.Range("E21").Value = '  +5.12%  '  # Volume(1h)
$ws.Range("D22").Value = '473.66'  # Price
$ws.Range("E22").Value = '  +1.67%  '  # Volume(1h)
$ws.Range("E23").Value = '  +1.08%  '  # Volume(1h)
$ws.Range("D24").Value = '0.0000164'  # Price
$ws.Range("E24").Value = '  +3.36%  '  # Volume(1h)
$ws.Range("D25").Value = '83.92'  # Price
$ws.Range("E25").Value = '  +0.91%  '  # Volume(1h)
$ws.Range("E26").Value = '  +3.74%  '  # Volume(1h)
$ws.Range("E27").Value = '  +1.37%  '  # Volume(1h)
$ws.Range("E28").Value = '  +5.77%  '  # Volume(1h)
$ws.Range("E29").Value = '  +0.20%  '  # Volume(1h)
$ws.Range("E30").Value = '  +1.76%  '  # Volume(1h)
$ws.Range("D31").Value = '4.025.62'  # Price
$ws.Range("E31").Value = '  +1.06%  '  # Volume(1h)
$ws.Range("E32").Value = '  +2.05%  '  # Volume(1h)
$ws.Range("D33").Value = '31.44'  # Price
$ws.Range("E33").Value = '  +2.13%  '  # Volume(1h)
$ws.Range("E34").Value = '  +1.51%  '  # Volume(1h)
$ws.Range("E35").Value = '  +1.83%  '  # Volume(1h)
$ws.Range("D36").Value = '3.837.55'  # Price
$ws.Range("E36").Value = '  +0.83%  '  # Volume(1h)
$ws.Range("D37").Value = '4.00'  # Price
$ws.Range("E37").Value = '  +23.57%  '  # Volume(1h)
$ws.Range("E38").Value = '  +1.44%  '  # Volume(1h)
$ws.Range("D39").Value = '1.03'  # Price
$ws.Range("E39").Value = '  +1.78%  '  # Volume(1h)
$ws.Range("E40").Value = '  +0.83%  '  # Volume(1h)
$ws.Range("D41").Value = '5.99'  # Price
$ws.Range("E41").Value = '  +2.51%  '  # Volume(1h)
$ws.Range("E42").Value = '  +0.06%  '  # Volume(1h)
$ws.Range("E43").Value = '  +3.37%  '  # Volume(1h)
$ws.Range("B44").Value = 'Stacks'  # Coin
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'  # Link
$ws.Range("D44").Value = '2.01'  # Price
$ws.Range("E44").Value = '  +2.26%  '  # Volume(1h)
$ws.Range("B45").Value = 'FLOKI'  # Coin
$ws.Range("C45").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'  # Link
$ws.Range("D45").Value = '0.000302'  # Price
$ws.Range("E45").Value = '  +14.26%  '  # Volume(1h)
$ws.Range("D46").Value = '424.19'  # Price
$ws.Range("E46").Value = '  +0.53%  '  # Volume(1h)
$ws.Range("E48").Value = '  +2.93%  '  # Volume(1h)
$ws.Range("D49").Value = '46.47'  # Price
$ws.Range("E49").Value = '  -1.44%  '  # Volume(1h)
$ws.Range("D50").Value = '142.61'  # Price
$ws.Range("E50").Value = '  -0.10%  '  # Volume(1h)
$ws.Range("D51").Value = '0.0360'  # Price
$ws.Range("E51").Value = '  +1.84%  '  # Volume(1h)

# Restore default cell style (NumberFormat change above would otherwise
# leave a stray style index applied to the cells).
$ws.Range("D2:E51").Style = "Normal"
